# Rename the worksheet from the generic default "Planilha1" to "Produtos"
# so the tab name matches the workbook's actual content (a product
# registration / "Cadastro Produtos" table) used in item 7 of the Power BI
# course. The "_xlnm._FilterDatabase" defined name (which references the
# sheet by name) is updated automatically by Excel when the sheet is renamed.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")
$ws.Name = "Produtos"
